# Daily attendance processing - 2026-01-28 21:46:47
# Swap the order of "Recorded By" entries in column G from "System, <email>"
# to "<email>, System" for the rows that currently hold that exact pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
